$d = $word.ActiveDocument

function Replace-ExactText($doc, [string]$oldText, [string]$newText) {
    $rng = $doc.Content
    $found = $rng.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found) {
        # Assign Range.Text directly (instead of Find.Execute's Replace:=wdReplaceOne)
        # so Word's smart-quote autocorrect does not mangle straight apostrophes.
        $rng.Text = $newText
    }
    return $found
}

# --- Section: "Drawing insights ... on preferred brands" ---
# Old Men/Women/Teens block is replaced (and reordered to Teens/Women/Men).

$oldMenBrands = "Men: Men often associate ice cream, particularly brands like Ezaki Glico and Ben & Jerry's, with comfort, cravings, and post-meal indulgence, sometimes going out of their way to find specific flavors. They value familiar brands and are influenced by presentation and variety of flavors when making purchase decisions.`n" + `
"Women: Women appreciate the convenience of readily available ice cream like Ezaki Glico, often purchasing it from convenience stores as a relaxing treat. They express disappointment when favorite flavors are discontinued and are sensitive to price, seeking online deals but finding some brands too expensive for regular purchase.`n" + `
"Teens: Teens consider ice cream a versatile treat, suitable for both individual snacks and sharing during social occasions, with brand preferences influenced by taste and value. They are open to switching brands for better flavors and appreciate brands that acknowledge portion control."

$newMenBrands = "Teens: Teens often associate ice cream with social occasions, sharing tubs with family or individual servings with friends. They are open to switching brands for better flavors and value freshness and natural taste.`n" + `
"Women: Women appreciate the convenience of readily available ice cream and may associate it with relaxation. Discontinued flavors can lead to disappointment, and they seek transparency from brands regarding these decisions.`n" + `
"Men: Men use ice cream as a comfort food and are drawn to brands with diverse flavors and convenient purchasing options. They are willing to seek out specific flavors they crave."

$found1 = Replace-ExactText $d $oldMenBrands $newMenBrands
Write-Host "Replace1: $found1"

# --- Section: "Drawing insights ... on frequency of purchase" ---

$oldFreq = "Men: Men purchase ice cream both less than once a month and more than once a week, influenced by factors like cravings, social context, and deals, with some prioritizing larger, economical tubs for frequent consumption.`n" + `
"Women: Women's ice cream purchase frequency varies from less than once a month to more than once a week, driven by factors like cravings, convenience, and variety, with some using it as a simple pleasure or reward.`n" + `
"Teens: Teens purchase ice cream from less than once a month to more than once a week, influenced by deals, social gatherings, and cravings, with some prioritizing unique flavors and smaller sizes for variety."

$newFreq = "Men: Men purchase ice cream both less than once a month and more than once a week, influenced by factors like deals, cravings, social context, and convenience, with brand preferences varying based on the occasion.`n" + `
"Women: Women's ice cream purchase frequency varies from less than once a month to more than once a week, driven by factors like cravings, convenience, variety, and brand loyalty, often viewing it as a simple pleasure or treat.`n" + `
"Teens: Teens purchase ice cream from less than once a month to more than once a week, influenced by deals, social gatherings, unique flavors, and convenience, with brand preferences varying based on quality and personal taste."

$found2 = Replace-ExactText $d $oldFreq $newFreq
Write-Host "Replace2: $found2"
